# Update "想去人数" (interested-count) figures in the generated sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 529
$ws1.Range("F5").Value = 241
$ws1.Range("F6").Value = 375
$ws1.Range("F7").Value = 234
$ws1.Range("F8").Value = 2262
$ws1.Range("F9").Value = 381
$ws1.Range("F10").Value = 5589
$ws1.Range("F11").Value = 131
$ws1.Range("F12").Value = 366

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 14

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 529
$ws4.Range("F6").Value = 241
$ws4.Range("F7").Value = 375
$ws4.Range("F8").Value = 234
$ws4.Range("F10").Value = 14
$ws4.Range("F11").Value = 2262
$ws4.Range("F12").Value = 381
$ws4.Range("F13").Value = 5589
$ws4.Range("F14").Value = 131
$ws4.Range("F15").Value = 366
